# Generate Report for Handoff
# Row 3 of every sheet refers to the "b.md" file. It has moved from
# "Handed back: in sync with en-US" to "Ready for handoff": a new handoff
# xliff was generated, and its handback is stale, so an error detail is
# now populated in the language sheets. Column P ("Error Detail") is widened
# to fit the new text.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e5a689d43004fe161f5135897bd37fc5cedf8b99/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1aac46d6e1d0ef13ff179ed469e6b66b7382a632/e2e/b.md."

# --- Overview sheet: row 3 is "b.md" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = "2016-08-16 06:33:17"

# --- zh-cn sheet: row 3 is "b.md" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 06:33:12"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is "b.md" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 06:33:17"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
